$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fix typo: "card.js" -> "cart.js" in A14 (the handleButtonClick row really belongs to cart.js)
$ws.Range("A14").Value = "cart.js"

# Update the line-number ranges for the cart.js functions (source file line numbers shifted)
$ws.Range("B12").Value = "128 à 189"
$ws.Range("B13").Value = "192 à 268"
$ws.Range("B14").Value = "283 à 444"
$ws.Range("B15").Value = "454 à 460"
$ws.Range("B16").Value = "470 à 476"
$ws.Range("B17").Value = "486 à 493"

# Reset the view: scroll back to top and select B4
$ws.Activate()
$ws.Range("B4").Select()
